$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 header-like numeric flags
$ws.Range("B2").Value = 1
$ws.Range("E2").Value = 0

# Row 4 data values
$ws.Range("B4").Value = 0.523252976771423
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.04650595354284603
$ws.Range("E4").Value = 0.7361990430222111
$ws.Range("F4").Value = -1
$ws.Range("G4").Value = 0.4723980860444221
$ws.Range("H4").Value = 0.7699568058175457
$ws.Range("J4").Value = 0.5399136116350913
